$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F is "dSF" -- update specific rows per repulled data / mean calculation
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 10
$ws.Range("F7").Value = -8
$ws.Range("F11").Value = -3
